# The presentation ships two theme parts:
#   ppt/theme/theme1.xml  -> bound to the (only) slide master, i.e. the
#                            theme actually used to render the slides.
#                            Originally the "Integral" / "Red Violet" theme.
#   ppt/theme/theme2.xml  -> bound only to the notes master ("Office Theme").
#
# The authored change swaps the two themes' content, so the slides end up
# using the plain "Office Theme" colours and the notes master ends up with
# the former "Integral" / "Red Violet" colours.
#
# The PowerPoint object model only exposes a mutable 12-slot theme colour
# scheme on objects that resolve back to the slide master's theme part
# (Slide.ThemeColorScheme / SlideRange.ThemeColorScheme / NotesPage.ThemeColorScheme
# all read & write the very same theme1.xml - there is no COM surface that
# reaches the notes-master-only theme2.xml). So we apply the achievable,
# faithful part of the edit: re-point the slide theme's 12 scheme colours
# (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) from the Integral/Red
# Violet palette to the Office palette that the diff moves onto theme1.xml.

$p = $ppt.ActivePresentation

function RGBVal($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

# Target palette (theme "Office"), in clrScheme order:
# dk1, lt1, dk2, lt2, accent1, accent2, accent3, accent4, accent5, accent6, hlink, folHlink
$officeColors = @(
    @(0x00, 0x00, 0x00),  # dk1
    @(0xFF, 0xFF, 0xFF),  # lt1
    @(0x44, 0x54, 0x6A),  # dk2
    @(0xE7, 0xE6, 0xE6),  # lt2
    @(0x5B, 0x9B, 0xD5),  # accent1
    @(0xED, 0x7D, 0x31),  # accent2
    @(0xA5, 0xA5, 0xA5),  # accent3
    @(0xFF, 0xC0, 0x00),  # accent4
    @(0x44, 0x72, 0xC4),  # accent5
    @(0x70, 0xAD, 0x47),  # accent6
    @(0x05, 0x63, 0xC1),  # hlink
    @(0x95, 0x4F, 0x72)   # folHlink
)

# Apply through the first slide's ThemeColorScheme - this writes straight
# into the shared slide-master theme part (ppt/theme/theme1.xml) that every
# slide/layout inherits its colours from.
$slide = $p.Slides.Item(1)
for ($i = 1; $i -le 12; $i++) {
    $c = $officeColors[$i - 1]
    $slide.ThemeColorScheme.Item($i).RGB = RGBVal $c[0] $c[1] $c[2]
}
